$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the "100 Home Project" background-processing row), which shifts
# the old row 3 up to become the new row 2.
$ws.Rows.Item(2).Delete()

# Leave the whole of (the new) row 2 selected, as happens after an entire-row delete.
$ws.Rows.Item(2).Select()
